$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A previously held plain usernames (kalesh220, swathi220, teja220, padma220).
# Replace them with full gmail addresses and turn them into mailto hyperlinks,
# matching the "Hyperlink" look already used in column B.
$ws.Range("A2").Value = "kalesh220@gmail.com"
$ws.Range("A3").Value = "pavan220@gmail.com"
$ws.Range("A4").Value = "teja220@gmail.com"
$ws.Range("A5").Value = "padmasri220@gmail.com"

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:kalesh220@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:pavan220@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:teja220@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:padmasri220@gmail.com")

# Match the Hyperlink cell style already used for column B.
$ws.Range("A2:A5").Style = "Hyperlink"

# Column A needs to be widened to fit the longer email addresses (target
# rendered width ~35.14 characters), and the "best fit" auto-sizing flag no
# longer applies since the width is now set explicitly.
$ws.Columns.Item(1).ColumnWidth = 34.33

# The active selection in the sheet view moved to I11.
$ws.Range("I11").Select()
